$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9963370561599731
$ws.Range("B1").Value = 1.832856059074402
$ws.Range("C1").Value = 5.587380886077881
$ws.Range("D1").Value = 1.60502827167511
$ws.Range("E1").Value = 0.6523777842521667
